$wb = $excel.ActiveWorkbook

# The "想去人数" (number of people interested) counts were updated for a
# handful of events. The same data is duplicated across the "展览" and
# "全部类型" worksheets, so apply the update to both.
$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2"  = 3011
    "F5"  = 6772
    "F6"  = 1764
    "F11" = 135
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
